$d = $word.ActiveDocument

# 1+3. "English" -> "ภาษาอังกฤษ" (both occurrences use same replacement text; both are
# matched fully inside their own run bounds with no differently-formatted predecessor
# at the exact match start, so this is safe as a single ReplaceAll)
$d.Content.Find.Execute("English", $true, $false, $false, $false, $false, $true, 1, $false, "ภาษาอังกฤษ", 2)

# 2. Language list line (leave the leading " / " literally in place, since that text's
# run immediately follows the hyperlink "English" run -- replacing starting at the very
# first character of that run would make it incorrectly inherit the hyperlink's format)
$d.Paragraphs(1).Range.Find.Execute("Portuguese / French / Thai / Vietnamese / Spanish", $true, $false, $false, $false, $false, $true, 1, $false, "ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน", 2)

# 4. "Brief" -> "บทย่อ"
$d.Content.Find.Execute("Brief", $true, $false, $false, $false, $false, $true, 1, $false, "บทย่อ", 2)

# 5. Brief description sentence
$d.Content.Find.Execute("An email sent to partners in the target country who have RSVPed no. It will be sent via customer.io", $true, $false, $false, $false, $false, $true, 1, $false, "An email sent to partners in the target country who have RSVPed no. โดยมันจะถูกส่งผ่านทาง customer.io", 2)

# 6. "Target audience" -> "กลุ่มเป้าหมาย"
$d.Content.Find.Execute("Target audience", $true, $false, $false, $false, $false, $true, 1, $false, "กลุ่มเป้าหมาย", 2)

# 7. "We'll miss you at the " -> "เราจะคิดถึงคุณจากที่ในงาน "
$d.Content.Find.Execute("We" + [char]8217 + "ll miss you at the ", $true, $false, $false, $false, $false, $true, 1, $false, "เราจะคิดถึงคุณจากที่ในงาน ", 2)

# 8. "Dear " -> "เรียนคุณ "
$d.Content.Find.Execute("Dear ", $true, $false, $false, $false, $false, $true, 1, $false, "เรียนคุณ ", 2)

# 9. ", " after [PARTNER NAME] -> " "
$d.Paragraphs(16).Range.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2)

# 10. "Thank you for taking the time..." -> Thai
$d.Content.Find.Execute("Thank you for taking the time to respond to our invitation to the upcoming ", $true, $false, $false, $false, $false, $true, 1, $false, "ขอขอบคุณที่สละเวลาตอบกลับคำเชิญของเราสำหรับงาน ", 2)

# 11. ". We were really looking forward to seeing you there." -> Thai
$d.Content.Find.Execute(". We were really looking forward to seeing you there.", $true, $false, $false, $false, $false, $true, 1, $false, " ที่กำลังจะมาถึง พวกเราหวังไว้ว่า จะได้พบเจอคุณที่นั่น", 2)

# 12. "Even though we're disappointed..." -> Thai
$d.Content.Find.Execute("Even though we" + [char]8217 + "re disappointed we can" + [char]8217 + "t meet you, we understand that scheduling conflicts and other commitments sometimes come up. ", $true, $false, $false, $false, $false, $true, 1, $false, "แม้ว่าเราจะผิดหวังที่ไม่สามารถพบคุณได้ แต่เราก็เข้าใจดีว่าปัญหาเกี่ยวกับตารางเวลาที่ขัดแย้งและภาระผูกพันอื่นๆ บางครั้งก็เกิดขึ้นได้ ", 2)

# 13. "If you're comfortable sharing it with us..." -> Thai
$d.Content.Find.Execute("If you" + [char]8217 + "re comfortable sharing it with us, we" + [char]8217 + "d like to know why you responded no. Please reply to this email as your feedback could help us make improvements in our event planning processes and better serve you in the future.", $true, $false, $false, $false, $false, $true, 1, $false, "หากคุณไม่ขัดข้องที่จะแบ่งปันกับเรา พวกเราก็ต้องการทราบว่า คุณตอบปฏิเสธคำเชิญเพราะอะไร โปรดตอบกลับอีเมล์นี้ เนื่องจากข้อคิดเห็นหรือคำติชมของคุณจะช่วยให้เราได้ปรับปรุงพัฒนากระบวนการวางแผนกิจกรรมของเราและให้บริการคุณได้ดียิ่งขึ้นในอนาคต", 2)

# 14. "We hope to see you at our future events. " -> Thai
$d.Content.Find.Execute("We hope to see you at our future events. ", $true, $false, $false, $false, $false, $true, 1, $false, "เราหวังว่า จะได้พบคุณในกิจกรรมของเราในอนาคต ", 2)

# 15. "If you have any questions, please contact us via " -> Thai
# This run is where a commentRangeStart marker sits, right at its very first character.
# Replacing a match that starts there pushes the marker to the end of the new text
# instead of leaving it at the start, so leave the leading "I" alone here and delete
# the stray leftover character afterwards (keeps commentRangeStart anchored correctly).
$d.Content.Find.Execute("f you have any questions, please contact us via ", $true, $false, $false, $false, $false, $true, 1, $false, "หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง ", 2)
$p22Start = $d.Paragraphs(22).Range.Start
$leftoverI = $d.Range($p22Start, $p22Start + 1)
$leftoverI.Text = ""

# 16. "live chat" -> "แชทสด" (hyperlink run). Leave the leading "l" alone in this
# replace (since matching the very first character of a hyperlink run loses its rPr),
# then delete the stray leftover "l" afterwards, scoped tightly to that hyperlink.
$d.Content.Find.Execute("ive chat", $true, $false, $false, $false, $false, $true, 1, $false, "แชทสด", 2)
$d.Hyperlinks(2).Range.Find.Execute("l", $true, $false, $false, $false, $false, $true, 1, $false, "", 1)

# 17. " or " (paragraph 22, after live chat hyperlink) -> " หรือทาง "
# Leave the leading space alone (it's the first character of the run right after the
# hyperlink, so replacing from its very start would lose/garble formatting); only
# replace "or " and put back the leading space unmodified.
$d.Paragraphs(22).Range.Find.Execute("or ", $true, $false, $false, $false, $false, $true, 1, $false, "หรือทาง ", 2)

# 18. ". " (paragraph 22, after WhatsApp hyperlink) -> " "
# Delete just the "." character (leaving the trailing space alone) using the
# WhatsApp hyperlink's end position, to avoid the same run-boundary formatting issue.
$whatsappLink = $d.Hyperlinks(3)
$dotRange = $d.Range($whatsappLink.Range.End, $whatsappLink.Range.End + 1)
$dotRange.Text = ""

# 19. "If you have any questions, please contact your country manager, " -> Thai
$d.Content.Find.Execute("If you have any questions, please contact your country manager, ", $true, $false, $false, $false, $false, $true, 1, $false, "หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ ", 2)

# 20. ", at " (paragraph 23, after [NAME]) -> " ที่ "
$d.Paragraphs(23).Range.Find.Execute(", at ", $true, $false, $false, $false, $false, $true, 1, $false, " ที่ ", 2)

# 21. " or " (paragraph 23, after [EMAIL ADDRESS]) -> " หรือ "
$d.Paragraphs(23).Range.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, " หรือ ", 2)

# 22. " (WhatsApp). " (paragraph 23, after [WHATSAPP NO]) -> " (WhatsApp) "
$d.Paragraphs(23).Range.Find.Execute(" (WhatsApp). ", $true, $false, $false, $false, $false, $true, 1, $false, " (WhatsApp) ", 2)

# 23. Comment text: "choose either one" -> "เลือกอย่างใดอย่างหนึ่ง"
# (Find.Execute on a comment's Range does not work reliably in this runtime since the
# comment's Range.Start/End are not relative to the comment's own story; assigning
# .Text directly on the comment Range works and is safe here because the comment's
# paragraph begins directly with this run, with no special run formatting to preserve.)
$d.Comments(1).Range.Text = "เลือกอย่างใดอย่างหนึ่ง"
